# TC31_Canine_Filter_Breed-Mixed_Neo4jData.xlsx
#
# The Neo4j "Stat" query was re-run (filter changed from breed 'Akita' to
# breed 'Mixed Breed'), which changes:
#   - StatOutput!A2:D2        -> the returned counts (files, samples, cases, study)
#   - StatOutput_Message!A18  -> the Cypher query text logged for that run
#
# The four numeric-looking counts must stay stored as TEXT (they were text
# before the edit too), so NumberFormat is forced to "@" before the value
# is written and cleared again afterwards so no visible formatting change
# is left behind.

$wb = $excel.ActiveWorkbook

# --- StatOutput: row 2 (file/sample/case/study counts for the new filter) ---
$wsStat = $wb.Worksheets.Item("StatOutput")

$statRange = $wsStat.Range("A2:D2")
$statRange.NumberFormat = "@"
$wsStat.Range("A2").Value = "82"
$wsStat.Range("B2").Value = "21"
$wsStat.Range("C2").Value = "28"
$wsStat.Range("D2").Value = "2"
$statRange.ClearFormats()

# --- StatOutput_Message: row 18 holds the Cypher query text that was run ---
$wsStatMsg = $wb.Worksheets.Item("StatOutput_Message")

$wsStatMsg.Range("A18").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Mixed Breed']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
